# Update column F ("dSF") values on Sheet1 per repull of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    2  = 1
    3  = -2
    4  = -5
    5  = 6
    6  = 5
    7  = 2
    8  = 1
    9  = -1
    10 = 1
    11 = 1
    12 = -3
    13 = 3
    14 = -3
    15 = 2
    16 = -5
    17 = -1
    18 = -2
    20 = -1
    21 = 6
    22 = 1
    23 = -5
    24 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
